# Remove duplicated teacher entries from the weekly schedule grid.
# (Commit message: "Listas sem duplicação de professores")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[-, -, 'MCT-3A-Processos de Usinagem 2', -]"

$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "[-, -, 'MCT-3A-Processos de Usinagem 2', -]"

$ws.Range("B4").Value = "-"
$ws.Range("D4").Value = "-"

$ws.Range("B6").Value = "[-, 'MCT-3A-Processos de Usinagem 2', -, -]"
$ws.Range("D6").Value = "-"

$ws.Range("B7").Value = "[-, 'MCT-3A-Processos de Usinagem 2', -, -]"
$ws.Range("D7").Value = "-"

$ws.Range("D14").Value = "-"

$ws.Range("D15").Value = "-"

$ws.Range("B18").Value = "-"
$ws.Range("F18").Value = "-"

$ws.Range("B19").Value = "-"
$ws.Range("F19").Value = "-"

$ws.Range("B20").Value = "-"
$ws.Range("D20").Value = "-"
$ws.Range("F20").Value = "-"

$ws.Range("F21").Value = "-"
